$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the PR reference hyperlink in B20 (same row as "Uthpal Suvarna").
# Passing TextToDisplay = Address first makes the engine record the hyperlink's
# "display" attribute as the target URL (matching real Excel's OOXML output);
# the cell's visible text is then set separately to the PR title so it reads
# as a normal hyperlink caption rather than the raw URL.
$prUrl = "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/3"
$prTitle = "Added Uthpal.Suvarna.txt by UthpalSuvarna " + [char]0x00B7 + " Pull Request #3 " + [char]0x00B7 + " dhavalkeerthi/MRIInterns2026A"

$ws.Hyperlinks.Add($ws.Range("B20"), $prUrl, "", "", $prUrl) | Out-Null
$ws.Range("B20").Value = $prTitle

# Restore the selection to the cell that was just edited, as the author would
# have left it positioned there when saving.
$ws.Range("B20").Select() | Out-Null
